# Added 2 scenarios for login feature
#
# The "NitroXLogin" sheet previously carried leftover Category/SearchText
# columns (copy/pasted from the "Amazon" sheet) and only had one real
# scenario row (row 2) pointed at the old Amazon hyperlink; rows 3 and 4
# were just placeholders ("NA"). This change removes the unrelated
# Category/SearchText columns and turns rows 3 and 4 into two more usable
# login scenarios by giving them the NitroX login URL (and drops the old
# hyperlink to amazon.in, replacing it with plain text of the new URL).
# It also makes NitroXLogin the active/selected sheet instead of NitroXHome.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NitroXLogin")

# Remove the leftover "Category" (E) and "SearchText" (F) columns - they
# don't apply to the login test data.
$ws.Range("E1:F4").EntireColumn.Delete()

# The old hyperlink (B2 -> https://www.amazon.in/) is no longer relevant;
# remove it before replacing the cell values below.
$ws.Hyperlinks.Delete()

# Populate the URL column for all three scenario rows with the NitroX
# login URL, adding 2 new usable scenarios (rows 3 & 4 used to be "NA").
$ws.Range("B2").Value = "https://test-nitrox.altono.app/"
$ws.Range("B3").Value = "https://test-nitrox.altono.app/"
$ws.Range("B4").Value = "https://test-nitrox.altono.app/"

# Resize column B now that it holds the (shorter) new URL text.
$ws.Columns("B").ColumnWidth = 25.75

# Make NitroXLogin the active sheet/selection (previously NitroXHome was).
[void]$ws.Activate()
[void]$ws.Range("F6").Select()
